# Limite em Empresa & Atualização de Saldo & Transação Cancelamento
#
# Insert 4 new error-code rows (21-24) into the "Transação de Combustível
# Veículo" error table on Planilha1, right after the existing row for
# error "55" (row 21) and before the blank separator row that precedes the
# "Transação de Combustível Equipamento" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Push rows 22 and below down by 4 to make room for the new entries.
$ws.Rows("22:25").Insert()

# New row 22: código 21
$ws.Cells.Item(22, 1).Value = "21"
$ws.Cells.Item(22, 2).Value = "Intervalo de Tempo inválido"
$ws.Cells.Item(22, 3).Value = "Tentativa de novo abastecimento realizada em um intervalo de tempo inferior ao mínimo permitido"

# New row 23: código 22
$ws.Cells.Item(23, 1).Value = "22"
$ws.Cells.Item(23, 2).Value = "Intervalo Percorrido inválido"
$ws.Cells.Item(23, 3).Value = "Tentativa de novo abastecimento realizada em um intervalo de kms percorridos inferior ao mínimo permitido"

# New row 24: código 23
$ws.Cells.Item(24, 1).Value = "23"
$ws.Cells.Item(24, 2).Value = "Horário inválido"
$ws.Cells.Item(24, 3).Value = "Tentativa de abastecimento fora do horário permitido"

# New row 25: código 24
$ws.Cells.Item(25, 1).Value = "24"
$ws.Cells.Item(25, 2).Value = "Qtde litros abastecidos inválida"
$ws.Cells.Item(25, 3).Value = "Tentativa de abastecimento acima da litragem máxima permitida"

# Widen column C to fit the new, longer descriptions (closest achievable
# value on the pixel-quantized column-width grid to the authored 88.23).
$ws.Columns("C").ColumnWidth = 87.33

# Restore the active selection recorded after the edit.
$ws.Range("C26").Select()
